$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 97
$ws.Range("I4").Value = 96.8
$ws.Range("J4").Value = 98
$ws.Range("K4").Value = 96.8
$ws.Range("L4").Value = 98
$ws.Range("M4").Value = 17.2
$ws.Range("N4").Value = -326

$ws.Range("H40").Value = 2093.524
$ws.Range("J40").Value = 2081.182
$ws.Range("L40").Value = 2081.182
$ws.Range("N40").Value = -2431.182

$ws.Range("H62").Value = 6961951.5
$ws.Range("I62").Value = 10113066
$ws.Range("J62").Value = 29500
$ws.Range("K62").Value = 10113066
$ws.Range("L62").Value = 29500
$ws.Range("M62").Value = -10112442
$ws.Range("N62").Value = -30748

$ws.Range("H65").Value = 6961951.5
$ws.Range("I65").Value = 10113066
$ws.Range("J65").Value = 29500
$ws.Range("K65").Value = 50565330
$ws.Range("L65").Value = 147500
$ws.Range("M65").Value = -50562210
$ws.Range("N65").Value = -153740

$ws.Range("H70").Value = 1667.5385
$ws.Range("I70").Value = 1075.4
$ws.Range("J70").Value = 2037.625
$ws.Range("K70").Value = 3226.2
$ws.Range("L70").Value = 6112.875
$ws.Range("M70").Value = -2956.2
$ws.Range("N70").Value = -6652.875

$ws.Range("H73").Value = 1667.5385
$ws.Range("I73").Value = 1075.4
$ws.Range("J73").Value = 2037.625
$ws.Range("K73").Value = 3226.2
$ws.Range("L73").Value = 6112.875
$ws.Range("M73").Value = -2290.2
$ws.Range("N73").Value = -7984.875

$ws.Range("H129").Value = 1043
$ws.Range("I129").Value = 533.3333
$ws.Range("J129").Value = 1083.2368
$ws.Range("K129").Value = 1599.9999
$ws.Range("L129").Value = 3249.7104
$ws.Range("M129").Value = 3400.0001
$ws.Range("N129").Value = -13249.7104

$ws.Range("H137").Value = 27028092
$ws.Range("I137").Value = 35715276
$ws.Range("J137").Value = 1299.7778
$ws.Range("K137").Value = 107145828
$ws.Range("L137").Value = 3899.3334
$ws.Range("M137").Value = -107143278
$ws.Range("N137").Value = -8999.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19004.72
$ws.Range("I32").Value = 3282.9524
$ws.Range("K32").Value = 3282.9524
$ws.Range("M32").Value = -2995.9524

$ws.Range("H61").Value = 4104.9375
$ws.Range("I61").Value = 3239.913
$ws.Range("J61").Value = 6315.5557
$ws.Range("K61").Value = 3239.913
$ws.Range("L61").Value = 6315.5557
$ws.Range("M61").Value = -3027.913
$ws.Range("N61").Value = -6739.5557

$ws.Range("H132").Value = 3619.6956
$ws.Range("I132").Value = 3028.647
$ws.Range("J132").Value = 5294.3335
$ws.Range("K132").Value = 9085.940999999999
$ws.Range("L132").Value = 15883.0005
$ws.Range("M132").Value = -6555.940999999999
$ws.Range("N132").Value = -20943.0005

$ws.Range("H136").Value = 4104.9375
$ws.Range("I136").Value = 3239.913
$ws.Range("J136").Value = 6315.5557
$ws.Range("K136").Value = 9719.739
$ws.Range("L136").Value = 18946.6671
$ws.Range("M136").Value = -7169.739
$ws.Range("N136").Value = -24046.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3998.6
$ws.Range("I134").Value = 3084.0557
$ws.Range("K134").Value = 9252.167099999999
$ws.Range("M134").Value = -6717.167099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6028.9165
$ws.Range("I31").Value = 2680
$ws.Range("J31").Value = 7316.9614
$ws.Range("K31").Value = 2680
$ws.Range("L31").Value = 7316.9614
$ws.Range("M31").Value = -2385
$ws.Range("N31").Value = -7906.9614

$ws.Range("H34").Value = 6028.9165
$ws.Range("I34").Value = 2680
$ws.Range("J34").Value = 7316.9614
$ws.Range("K34").Value = 2680
$ws.Range("L34").Value = 7316.9614
$ws.Range("M34").Value = -2478
$ws.Range("N34").Value = -7720.9614

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H122").Value = 887.5714
$ws.Range("I122").Value = 887.5714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2662.7142
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -212.7142000000003
$ws.Range("N122").ClearContents()

$ws.Range("H131").Value = 35000
$ws.Range("J131").Value = 35000
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1592.9524
$ws.Range("J34").Value = 1513.0526
$ws.Range("L34").Value = 4539.1578
$ws.Range("N34").Value = -4707.1578

$ws.Range("H39").Value = 8119.231
$ws.Range("J39").Value = 8119.231
$ws.Range("L39").Value = 24357.693
$ws.Range("N39").Value = -24945.693

$ws.Range("H55").Value = 2375
$ws.Range("J55").Value = 2833.3333
$ws.Range("L55").Value = 8499.999899999999
$ws.Range("N55").Value = -8853.999899999999

$ws.Range("H75").Value = 752
$ws.Range("I75").Value = 186.33333
$ws.Range("J75").Value = 1176.25
$ws.Range("K75").Value = 558.99999
$ws.Range("L75").Value = 3528.75
$ws.Range("M75").Value = 439.00001
$ws.Range("N75").Value = -5524.75

$ws.Range("H76").Value = 6000
$ws.Range("J76").Value = 6000
$ws.Range("L76").Value = 18000
$ws.Range("N76").Value = -18766

$ws.Range("H78").Value = 752
$ws.Range("I78").Value = 186.33333
$ws.Range("J78").Value = 1176.25
$ws.Range("K78").Value = 1676.99997
$ws.Range("L78").Value = 10586.25
$ws.Range("M78").Value = 3315.00003
$ws.Range("N78").Value = -20570.25

$ws.Range("H79").Value = 6000
$ws.Range("J79").Value = 6000
$ws.Range("L79").Value = 18000
$ws.Range("N79").Value = -20652

$ws.Range("H80").Value = 1213.2858
$ws.Range("J80").Value = 1332.1666
$ws.Range("L80").Value = 3996.4998
$ws.Range("N80").Value = -5868.4998

$ws.Range("H83").Value = 1213.2858
$ws.Range("J83").Value = 1332.1666
$ws.Range("L83").Value = 11989.4994
$ws.Range("N83").Value = -21349.4994

$ws.Range("H126").Value = 3570.2856
$ws.Range("J126").Value = 4586.4
$ws.Range("L126").Value = 13759.2
$ws.Range("N126").Value = -23639.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5689.5625
$ws.Range("I70").Value = 6192.9
$ws.Range("J70").Value = 4850.6665
$ws.Range("K70").Value = 6192.9
$ws.Range("L70").Value = 4850.6665
$ws.Range("M70").Value = -5922.9
$ws.Range("N70").Value = -5390.6665

$ws.Range("H73").Value = 5689.5625
$ws.Range("I73").Value = 6192.9
$ws.Range("J73").Value = 4850.6665
$ws.Range("K73").Value = 6192.9
$ws.Range("L73").Value = 4850.6665
$ws.Range("M73").Value = -5256.9
$ws.Range("N73").Value = -6722.6665

$ws.Range("H122").Value = 2277.4443
$ws.Range("I122").Value = 1999.8334
$ws.Range("J122").Value = 2832.6667
$ws.Range("K122").Value = 5999.5002
$ws.Range("L122").Value = 8498.000100000001
$ws.Range("M122").Value = -3549.5002
$ws.Range("N122").Value = -13398.0001

$ws.Range("H123").Value = 15949.429
$ws.Range("J123").Value = 15949.429
$ws.Range("L123").Value = 15949.429
$ws.Range("N123").Value = -20849.429

$ws.Range("H132").Value = 2620.037
$ws.Range("I132").Value = 2492.5652
$ws.Range("J132").Value = 3353
$ws.Range("K132").Value = 7477.6956
$ws.Range("L132").Value = 10059
$ws.Range("M132").Value = -4947.6956
$ws.Range("N132").Value = -15119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4546516
$ws.Range("I16").Value = 6667383.5
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 6667383.5
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -6667213.5
$ws.Range("N16").Value = -2140

$ws.Range("H40").Value = 3185.2
$ws.Range("I40").Value = 1734.6666
$ws.Range("J40").Value = 3441.1765
$ws.Range("K40").Value = 1734.6666
$ws.Range("L40").Value = 3441.1765
$ws.Range("M40").Value = -1598.6666
$ws.Range("N40").Value = -3713.1765

$ws.Range("H122").Value = 3729.5
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3767.8948
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 11303.6844
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -16203.6844

$ws.Range("H136").Value = 6949.5
$ws.Range("I136").Value = 6949.5
$ws.Range("K136").Value = 20848.5
$ws.Range("M136").Value = -18298.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 326.63635
$ws.Range("I107").Value = 261.625
$ws.Range("K107").Value = 784.875
$ws.Range("M107").Value = 1135.125

$ws.Range("H136").Value = 3249.1667
$ws.Range("I136").Value = 1998.75
$ws.Range("K136").Value = 5996.25
$ws.Range("M136").Value = -3446.25

